# Applies the e1 report edit:
#  1. Update the "Date and Time Generated" timestamp string.
#  2. Clear the data row (row 16) contents, leaving formatting/style intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update generated timestamp (cell A9)
$ws.Range("A9").Value = "Date and Time Generated: July 20, 2024 05:58:38 PM"

# 2. Clear the contents of row 16 (columns A:AE had values; AF was already empty)
#    Use ClearContents so number formatting / styles on the cells are preserved.
$ws.Range("A16:AE16").ClearContents()
